$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.952.92"
$ws.Range("E2").Value = "  -3.18%  "

$ws.Range("D3").Value = "1.681.65"
$ws.Range("E3").Value = "  -2.79%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.003"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.59%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "309.55"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.15%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.9980"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.94%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.3671"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -2.69%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3375"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -5.80%  "

$ws.Range("E9").Value = "  -5.38%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "1.177"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -3.78%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07332"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -2.65%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.9983"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +1.07%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "6.199"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -2.49%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "20.54"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -5.36%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "6.847"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -2.37%  "

$ws.Range("D16").Value = "1.678.10"
$ws.Range("E16").Value = "  -3.08%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.00001102"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -3.71%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.06611"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -2.26%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.9970"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.01%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "82.31"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -4.09%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "16.88"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -2.91%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.200"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -2.38%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "12.65"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.97%  "

$ws.Range("D24").Value = "24.874.64"
$ws.Range("E24").Value = "  -3.26%  "

$ws.Range("E25").Value = "  -0.26%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.703"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -4.45%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "150.64"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -3.02%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "19.83"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -2.70%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.285"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +9.93%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "130.45"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -2.50%  "

$ws.Range("D31").Value = "1.869.28"
$ws.Range("E31").Value = "  -2.71%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "6.513"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -2.88%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "4.158"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +1.49%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "13.51"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +1.06%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.08615"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.54%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.734"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -2.55%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "5.447"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -1.97%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.06484"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -3.00%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.02347"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -4.74%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "8.756"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -4.48%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.2170"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.73%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.246"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -2.75%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.6284"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -2.24%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.9968"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.93%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "13.47"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -1.50%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "3.794"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -2.09%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.5987"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -3.87%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.043"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -3.90%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "125.75"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -4.48%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.07171"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -3.82%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "77.52"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -1.54%  "
